$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 (I0) and J1 (IF), formatted like the existing header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data for the new I0 / IF columns, rows 2-71
$i0Values = @(
    6, 6, 3, 7, 7, 10, 8, 5, 7, 7, 7, 6, 8, 8, 8, 8, 7, 8, 8, 9, 6, 7, 8, 8, 5, 7, 6, 8, 6, 7, 6, 7, 7, 5, 8, 7, 6, 8, 8, 8, 10, 7, 6, 9, 10, 8, 7, 8, 9, 7, 7, 8, 8, 7, 7, 8, 7, 6, 8, 7, 7, 6, 6, 6, 8, 9, 6, 5, 5, 5
)
$ifValues = @(
    7, 7, 3, 7, 7, 10, 8, 5, 7, 7, 7, 6, 8, 9, 8, 8, 8, 8, 8, 9, 6, 7, 8, 8, 6, 7, 6, 8, 6, 7, 6, 7, 8, 5, 8, 7, 6, 8, 8, 8, 10, 8, 6, 9, 10, 8, 7, 8, 9, 7, 7, 8, 8, 8, 7, 8, 8, 6, 8, 7, 8, 6, 6, 6, 8, 9, 6, 5, 5, 5
)

for ($idx = 0; $idx -lt $i0Values.Length; $idx++) {
    $r = $idx + 2
    $ws.Cells.Item($r, 9).Value = $i0Values[$idx]
    $ws.Cells.Item($r, 10).Value = $ifValues[$idx]
}
